$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header rename
$ws.Range("G1").Value = "S Tag"

# Speaker column (D) renames per row
$speakerUpdates = @(
    @(2, "T"), @(3, "T"), @(4, "T"), @(5, "T"), @(6, "T"), @(7, "T"), @(8, "T"), @(9, "T"), @(10, "T"), @(11, "T"),
    @(12, "T"), @(13, "S"), @(14, "T"), @(15, "T"), @(16, "S"), @(17, "T"), @(18, "S"), @(19, "T"), @(20, "T"),
    @(21, "T"), @(22, "T"), @(23, "T"), @(24, "T"), @(25, "T"), @(26, "T"), @(27, "T"), @(28, "S"), @(29, "T2"),
    @(30, "S"), @(31, "T2"), @(32, "S"), @(33, "T2"), @(34, "S"), @(35, "T2"), @(36, "S"), @(37, "T2"), @(38, "S"),
    @(39, "T2"), @(40, "T2"), @(41, "T2"), @(42, "T2"), @(43, "SS"), @(44, "T"), @(45, "T"), @(46, "S"), @(47, "T"),
    @(48, "T"), @(49, "T"), @(50, "T"), @(51, "T"), @(52, "S"), @(53, "T"), @(54, "T"), @(55, "T"), @(56, "S"),
    @(57, "T"), @(58, "T"), @(59, "T"), @(60, "T"), @(61, "T")
)

foreach ($update in $speakerUpdates) {
    $row = $update[0]
    $value = $update[1]
    $ws.Cells.Item($row, 4).Value = $value
}

# Teacher Tag (F) text updates
$ws.Range("F35").Value = "3 - Getting Ss to Relate"
$ws.Range("F37").Value = "3 - Getting Ss to Relate"

# Student Tag (G) text update
$ws.Range("G36").Value = "2 - Relating to Another S"
